# Add the new worksheet "Tabelle2" right after "Tabelle1" (becomes the active sheet,
# matching activeTab="1" on the workbook and tabSelected="1" on the new sheet).
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "Tabelle2"

# --- Copy over just the cell formatting (styles/borders/number formats) from the
# --- corresponding ranges on Tabelle1, cell-range by cell-range, so that no
# --- superfluous blank cells/rows are introduced and no new style entries are minted.
$ws1.Range("B2:F2").Copy()
$ws2.Range("B2").PasteSpecial(-4122)

$ws1.Range("H5").Copy()
$ws2.Range("H5").PasteSpecial(-4122)

$ws1.Range("H6:H9").Copy()
$ws2.Range("H6").PasteSpecial(-4122)

$ws1.Range("H10").Copy()
$ws2.Range("H10").PasteSpecial(-4122)

$ws1.Range("I6:I9").Copy()
$ws2.Range("I6").PasteSpecial(-4122)

$ws1.Range("I10").Copy()
$ws2.Range("I10").PasteSpecial(-4122)

$ws1.Range("J6:J9").Copy()
$ws2.Range("J6").PasteSpecial(-4122)

$ws1.Range("J10").Copy()
$ws2.Range("J10").PasteSpecial(-4122)

$ws1.Range("C10:G10").Copy()
$ws2.Range("C10").PasteSpecial(-4122)

$ws1.Range("H11").Copy()
$ws2.Range("H11").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# --- Header / label text (reuses the existing shared strings) ---
$ws2.Range("B2").Value = "Evaluation"
$ws2.Range("C2").Value = "fixed 4 6-Tupels 100k TDNT2 afterState.agt.zip"
$ws2.Range("C5").Value = "highest tile"
$ws2.Range("H5").Value = "sum"
$ws2.Range("I5").Value = "percent"
$ws2.Range("J5").Value = "percent cum"
$ws2.Range("C6").Value = "< 1024"

# --- Row labels (highest tile values) ---
$ws2.Range("C7").Value = 1024
$ws2.Range("C8").Value = 2048
$ws2.Range("C9").Value = 4096
$ws2.Range("C10").Value = 8192

# --- Raw sample data (columns D:G, rows 6:10) ---
$ws2.Range("D6").Value = 2
$ws2.Range("E6").Value = 0
$ws2.Range("F6").Value = 0
$ws2.Range("G6").Value = 0

$ws2.Range("D7").Value = 4
$ws2.Range("E7").Value = 1
$ws2.Range("F7").Value = 3
$ws2.Range("G7").Value = 1

$ws2.Range("D8").Value = 8
$ws2.Range("E8").Value = 13
$ws2.Range("F8").Value = 7
$ws2.Range("G8").Value = 10

$ws2.Range("D9").Value = 22
$ws2.Range("E9").Value = 25
$ws2.Range("F9").Value = 31
$ws2.Range("G9").Value = 22

$ws2.Range("D10").Value = 14
$ws2.Range("E10").Value = 11
$ws2.Range("F10").Value = 9
$ws2.Range("G10").Value = 17

# --- Formulas; whole-range assignment produces shared formula groups just like
# --- the ones already present on Tabelle1. The order of creation below matches
# --- the row order on Tabelle1 so the shared-group indices (si="0","1","2") line
# --- up the same way: J6:J8, then H7:H10, then E11:H11.
$ws2.Range("H6").Formula = "=SUM(D6:G6)"

$ws2.Range("I6").Formula = "=H6/H`$11"
$ws2.Range("I7").Formula = "=H7/H`$11"
$ws2.Range("I8").Formula = "=H8/H`$11"
$ws2.Range("I9").Formula = "=H9/H`$11"
$ws2.Range("I10").Formula = "=H10/H`$11"

$ws2.Range("J6:J8").Formula = "=J7+I6"

$ws2.Range("H7:H10").Formula = "=SUM(D7:G7)"

$ws2.Range("J9").Formula = "=J10+I9"
$ws2.Range("J10").Formula = "=I10"

$ws2.Range("D11").Formula = "=SUM(D6:D10)"
$ws2.Range("E11:H11").Formula = "=SUM(E6:E10)"

# --- Page setup (matches Tabelle1: A4-ish custom margins, portrait, paper size 9) ---
$ws2.PageSetup.PaperSize = 9
$ws2.PageSetup.Orientation = 1
$ws2.PageSetup.TopMargin = 56.6929134
$ws2.PageSetup.BottomMargin = 56.6929134
$ws2.PageSetup.LeftMargin = 50.4
$ws2.PageSetup.RightMargin = 50.4
$ws2.PageSetup.HeaderMargin = 21.6
$ws2.PageSetup.FooterMargin = 21.6

# --- Selection state for the new sheet (matches the captured sheetView) ---
$ws2.Range("D15").Select()

Write-Output "done"
